$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.281343579292297
$ws.Range("B1").Value = 1.764230847358704
$ws.Range("C1").Value = 1.641840815544128
$ws.Range("D1").Value = 4.960494041442871
$ws.Range("E1").Value = 1.38152015209198
